$d = $word.ActiveDocument

# Namespace declaration used for the inserted OOXML fragments.
$wns = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1. Remove the _GoBack bookmark from the "(pag 38...)" paragraph ---
# It is re-created later on, further down in the document, inside the new
# "@Autowired" paragraph. Deleting the bookmark object (rather than
# rewriting the whole paragraph) keeps every other attribute of that
# paragraph/run untouched.
$d.Bookmarks.Item("_GoBack").Delete()

# --- 2. Replace the "ApplicationContext es una interfaz..." paragraph with
#        the new block of paragraphs from the edit. ---
$count = $d.Paragraphs.Count
$appCtxIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -like "ApplicationContext es una interfaz*") {
        $appCtxIndex = $i
        break
    }
}

$appCtxPara = $d.Paragraphs.Item($appCtxIndex)
$newBlock = @'
<w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:t>ApplicationContext es una interfaz de Spring que permite instanciar de manera transparente los “beans” que hayan sido declarados, bien en el xml de configuración (app-context.xml), bien mediante anotaciones.</w:t>
      </w:r>
      <w:r>
        <w:t>s</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Dependency pull:</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> por ejemplo, JNDI, utilizado para buscar componentes EJB en EJB 2.1 o anteriores</w:t>
      </w:r>
      <w:r>
        <w:t>.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:t>El propio Spring utiliza ese mecanismo cuando hace un ctx.getBean(&lt;nombre&gt;, &lt;clase&gt;);</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve">Contextualized dependency lookup: </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">es similar al modelo “dependency pull”, pero </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">en lugar de acceder a un registro central, la búsqueda se hace directamente sobre el contenedor que gestiona el objeto. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve">Inyección de dependencias por constructor: </w:t>
      </w:r>
      <w:r>
        <w:t>obvio</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t xml:space="preserve">Inyección por métodos “setter”: </w:t>
      </w:r>
      <w:r>
        <w:t>es el mecanismo más utilizado de IoC</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">, y de los más sencillos. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:t>Comparación entre los “Injection” y “Lookup”:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
        <w:rPr>
          <w:b/>
        </w:rPr>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">En general, es mucho mejor el mecanismo de inyección, ya que, entre otras cosas, tiene un impacto cero en los componentes que lo utilizan. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">El componente central del contenedor de inyección de dependencias en Spring es la interfaz BeanFactory. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">En Spring, un “bean” es cualquier componente que sea gestionado por el contenedor. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:t>La configuración de la BeanFactory se lleva a cabo, internamente, mediante las clases que implementan la interfaz BeanDefinition</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> (una implementación es, por ejemplo, DefaultListableBeanFactory)</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">Externamente esa configuración se hace con algún fichero .xml (ahora ya con anotaciones). </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
        </w:rPr>
        <w:lastRenderedPageBreak/>
        <w:t>ApplicationContext</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">: </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">interfaz que </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">hereda de BeanFactory. </w:t>
      </w:r>
      <w:r>
        <w:t>Tiene</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> en sus implementaciones</w:t>
      </w:r>
      <w:r>
        <w:t>, además de los servicios de DI, servicios de A</w:t>
      </w:r>
      <w:r>
        <w:t>OP,</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>internacionalización, manejo de eventos y otros</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t xml:space="preserve">Configuración básica: </w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:t>Fichero app-context-annotation.xml</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:pStyle w:val="Standard"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">@Autowired: tiene un equivalente </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">(no del todo. Soporta el parámetro “name) </w:t>
      </w:r>
      <w:r>
        <w:t>en el estándar de Java que es @Resource(</w:t>
      </w:r>
      <w:r>
        <w:t>&lt;id&gt;</w:t>
      </w:r>
      <w:r>
        <w:t>)</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">. Otro equivalente </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">(éste sí) </w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
      <w:r>
        <w:t xml:space="preserve">de JEE es @Inject, también soportado por Spring. </w:t>
      </w:r>
    </w:p>
'@
$newBlock = $newBlock.Replace("<w:p>", "<w:p " + $wns + ">")
$appCtxPara.Range.InsertXML($newBlock)
